$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '33.885.41'
$ws.Range('E2').Value = '  +10.40%  '

$ws.Range('D3').Value = '1.809.56'
$ws.Range('E3').Value = '  +7.36%  '

$ws.Range('E4').Value = '  +0.05%  '

$ws.Range('D5').Value = "'227.83"
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +3.50%  '

$ws.Range('D6').Value = "'0.540"
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +3.16%  '

$ws.Range('D7').Value = "'0.999"
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +0.00%  '

$ws.Range('D8').Value = "'30.96"
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +1.84%  '

$ws.Range('D9').Value = "'47.31"
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +6.69%  '

$ws.Range('E10').Value = '  +5.47%  '

$ws.Range('D11').Value = "'0.0663"
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +6.03%  '

$ws.Range('D12').Value = "'0.0929"
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +2.47%  '

$ws.Range('D13').Value = '2.070.33'
$ws.Range('E13').Value = '  +7.38%  '

$ws.Range('D14').Value = '1.811.61'
$ws.Range('E14').Value = '  +7.60%  '

$ws.Range('E15').Value = '  +2.51%  '

$ws.Range('D16').Value = '33.893.96'
$ws.Range('E16').Value = '  +10.41%  '

$ws.Range('D17').Value = "'10.05"
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -3.87%  '

$ws.Range('E18').Value = '  +6.35%  '

$ws.Range('D19').Value = "'69.17"
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +4.38%  '

$ws.Range('D20').Value = "'255.10"
$ws.Range('D20').Style = "Normal"

$ws.Range('E21').Value = '  +3.54%  '

$ws.Range('E22').Value = '  +0.03%  '

$ws.Range('D23').Value = "'10.40"
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +1.91%  '

$ws.Range('E24').Value = '  +0.37%  '

$ws.Range('E25').Value = '  +1.03%  '

$ws.Range('D26').Value = "'158.87"
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +0.76%  '

$ws.Range('D27').Value = "'16.40"
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +3.62%  '

$ws.Range('E28').Value = '  +2.89%  '

$ws.Range('E29').Value = '  +5.02%  '

$ws.Range('E30').Value = '  +0.04%  '

$ws.Range('D31').Value = "'3.81"
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +9.59%  '

$ws.Range('D32').Value = "'0.0508"
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +1.92%  '

$ws.Range('E33').Value = '  +5.35%  '

$ws.Range('E34').Value = '  +6.45%  '

$ws.Range('D35').Value = '1.537.21'
$ws.Range('E35').Value = '  +1.79%  '

$ws.Range('E36').Value = '  +2.89%  '

$ws.Range('B37').Value = 'MinaProtocolToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/J7st_qGwz+minaprotocoltoken-mina'
$ws.Range('D37').Value = "'1.43"
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +251.46%  '

$ws.Range('B38').Value = 'TrustWalletToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D38').Value = "'1.07"
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +4.03%  '

$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').Value = "'0.0186"
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +4.30%  '

$ws.Range('B40').Value = 'Aave'
$ws.Range('C40').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D40').Value = "'83.47"
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -1.30%  '

$ws.Range('B41').Value = 'ImmutableX'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D41').Value = "'0.616"
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +5.33%  '

$ws.Range('B42').Value = 'MXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D42').Value = "'2.82"
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +3.92%  '

$ws.Range('B43').Value = 'HuobiToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D43').Value = "'2.33"
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +0.22%  '

$ws.Range('B44').Value = 'ARBITRUM'
$ws.Range('C44').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D44').Value = "'0.901"
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +7.74%  '

$ws.Range('B45').Value = 'RenderToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D45').Value = "'2.10"
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +5.43%  '

$ws.Range('B46').Value = 'Kaspa'
$ws.Range('C46').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D46').Value = "'0.0520"
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +4.10%  '

$ws.Range('B47').Value = 'WEMIXToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D47').Value = "'1.07"
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +4.40%  '

$ws.Range('B48').Value = 'RocketPoolETH'
$ws.Range('C48').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D48').Value = '1.957.95'
$ws.Range('E48').Value = '  +7.37%  '

$ws.Range('B49').Value = 'PaxDollar'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D49').Value = "'0.999"
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -0.08%  '

$ws.Range('D50').Value = "'5.63"
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +3.98%  '

$ws.Range('D51').Value = "'51.95"
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +0.14%  '
